$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a text value to a cell while preventing Excel from
# auto-converting date/number-looking strings (e.g. "2026-02-17") into a
# serial date. We temporarily mark the cell as Text, assign the value, then
# restore the cell style to "Normal" so no extra formatting is left behind.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ===========================================================================
# Sheet "Summary": refresh aggregate metrics
# ===========================================================================
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.83   # Current Capital
$summary.Range("B4").Value = -1.18    # Total P&L $
$summary.Range("B6").Value = 119      # Total Trades
$summary.Range("B8").Value = 60       # Losing Trades
$summary.Range("B9").Value = 36.13    # Win Rate %

# ===========================================================================
# Sheet "Strategy Status": refresh MarketMaking row (row 4)
# ===========================================================================
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.83     # Capital
$status.Range("D4").Value = 119       # Trades
$status.Range("E4").Value = -1.18     # P&L $
$status.Range("F4").Value = -1.17     # P&L %
$status.Range("G4").Value = 36.13     # Win Rate %

# ===========================================================================
# New closed trade #119 appended as row 120 on both the "All Trades" sheet
# and the "MarketMaking" sheet (they mirror the same trade log).
# ===========================================================================
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A120").Value = 119

    Set-TextValue $ws.Range("B120") "2026-02-17"
    Set-TextValue $ws.Range("C120") "16:03:50"
    Set-TextValue $ws.Range("D120") "MarketMaking"
    Set-TextValue $ws.Range("E120") "UP"

    $ws.Range("F120").Value = 0.91
    $ws.Range("G120").Value = 0.9

    Set-TextValue $ws.Range("H120") "CLOSED"

    $ws.Range("I120").Value = -1.0989
    $ws.Range("J120").Value = -0.01
    $ws.Range("K120").Value = 98.83
    $ws.Range("L120").Value = 0
    $ws.Range("M120").Value = 0
    $ws.Range("N120").Value = 0.6

    Set-TextValue $ws.Range("O120") "Normal spread capture: 19600 bps"
    Set-TextValue $ws.Range("P120") "early_exit"

    $ws.Range("Q120").Value = 0.11
}
